# "Generate Report for Handback" - refresh the localization-status report:
#   * Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   * Handback timestamps for zh-cn / de-de are refreshed
#   * The (now resolved) "handback file is not latest" error detail is cleared
#   * A few report columns are resized to fit the new content

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$newStatus = "Handed back: in sync with en-US"

$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws2.Range("C2").Value = $newStatus
$ws3.Range("C2").Value = $newStatus

# --- Refresh "Latest Handback DateTime" for each locale ---
$ws2.Range("K2").Value = "2016-08-16 22:48:28"
$ws3.Range("K2").Value = "2016-08-16 22:48:35"

# --- Clear the stale "Error Detail" message now that handback is in sync ---
$ws2.Range("P2").Value = ""
$ws3.Range("P2").Value = ""

# --- Resize columns to fit the refreshed content ---
# (ColumnWidth is expressed in characters; Excel snaps the stored width to
# whole-pixel boundaries, so we pick the character width that lands on the
# nearest achievable pixel boundary to the target column width.)
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(16).ColumnWidth = 12.833333333333334

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(16).ColumnWidth = 12.833333333333334
